$d = $word.ActiveDocument

# --- Paragraph 1: "Option to set a destination ..." -> highlight green ---
$target1 = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*Option to set a destination to the directory patient folders will be created*") {
        $target1 = $p
        break
    }
}
if ($target1 -ne $null) {
    $target1.Range.Font.HighlightColorIndex = "wdBrightGreen"
}

# --- Paragraph 2: "(ME): set this via a file dialog??... PRIO:2)" -> magenta to green ---
$target2 = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*set this via a file dialog*PRIO:2*") {
        $target2 = $p
        break
    }
}
if ($target2 -ne $null) {
    $target2.Range.Font.HighlightColorIndex = "wdBrightGreen"
}

# The hyperlink run inside that paragraph doesn't always pick up the
# highlight change from the enclosing paragraph range, so touch it directly.
for ($i = 1; $i -le $d.Hyperlinks.Count; $i++) {
    $h = $d.Hyperlinks.Item($i)
    if ($h.Range.Text -eq "https://github.com/dfranx/ImFileDialog") {
        $h.Range.Font.HighlightColorIndex = "wdBrightGreen"
    }
}

Write-Host "Highlight updated to green for the two target paragraphs."
